$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 162, shifting all
# subsequent rows (162-238) down to (164-240).
$ws.Rows.Item(162).Insert()
$ws.Rows.Item(162).Insert()

# New row 162: Albahaca, Primera, Provincia de Chacabuco
$ws.Cells.Item(162, 1).Value = 9
$ws.Cells.Item(162, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(162, 3).Value = "Metropolitana"
$ws.Cells.Item(162, 4).Value = 44523
$ws.Cells.Item(162, 5).Value = 13
$ws.Cells.Item(162, 6).Value = 100112052
$ws.Cells.Item(162, 7).Value = "Albahaca"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 52
$ws.Cells.Item(162, 11).Value = 6000
$ws.Cells.Item(162, 12).Value = 7000
$ws.Cells.Item(162, 13).Value = 6500
$ws.Cells.Item(162, 14).Value = "$/docena de matas"
$ws.Cells.Item(162, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(162, 16).Value = 1083
$ws.Cells.Item(162, 17).Value = 6
$ws.Cells.Item(162, 18).Value = "Hortaliza"

# New row 163: Albahaca, Primera, Provincia del Elquí
$ws.Cells.Item(163, 1).Value = 9
$ws.Cells.Item(163, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(163, 3).Value = "Metropolitana"
$ws.Cells.Item(163, 4).Value = 44523
$ws.Cells.Item(163, 5).Value = 13
$ws.Cells.Item(163, 6).Value = 100112052
$ws.Cells.Item(163, 7).Value = "Albahaca"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 34
$ws.Cells.Item(163, 11).Value = 7000
$ws.Cells.Item(163, 12).Value = 8000
$ws.Cells.Item(163, 13).Value = 7500
$ws.Cells.Item(163, 14).Value = "$/docena de matas"
$ws.Cells.Item(163, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(163, 16).Value = 1250
$ws.Cells.Item(163, 17).Value = 6
$ws.Cells.Item(163, 18).Value = "Hortaliza"
